$d = $word.ActiveDocument

# --- 1) "Entered seed URLs must be full urls. EG: " paragraph: no text change,
#         just normalize the run that had the spell-check split around "urls".
$t = "Entered seed URLs must be full urls. EG: "
$d.Content.Find.Execute($t, $false, $false, $false, $false, $false, $true, 1, $false, $t, 2)

# --- 2) "You will then be asked to enter a keyword to search for. Do that."
#         -> "...keyword to search for in each pages meta keywords. Do that."
$d.Content.Find.Execute(
    "keyword to search for. Do that.", $false, $false, $false, $false, $false,
    $true, 1, $false,
    "keyword to search for in each pages meta keywords. Do that.", 2)

# --- 3) "SpiderLeg" heading: merge "SpiderLeg" + ":" runs, dropping the spell-check split.
$t = "SpiderLeg:"
$d.Content.Find.Execute($t, $false, $false, $false, $false, $false, $true, 1, $false, $t, 2)

# --- 4) getHyperLink paragraph: normalize away the spell-check splits around
#         "getHyperLink", "urls" and "Youtube" (text itself is unchanged).
$t = "The getHyperLink method also validates any links it finds before returning them. It checks to ensure they a"
$d.Content.Find.Execute($t, $false, $false, $false, $false, $false, $true, 1, $false, $t, 2)

$t = "re valid urls, and ensures they are not on the exclusion list. The exclusion list includes large, complex sites such as Twitter, Facebook, and Youtube. It also excludes pages in formats that it cannot scrape, such as .asp sites."
$d.Content.Find.Execute($t, $false, $false, $false, $false, $false, $true, 1, $false, $t, 2)

# --- 5) DFS -> BFS
$d.Content.Find.Execute(
    "This class does the actual crawling. It is given a list of seed URLs, and effectively conducts a DFS using the seeds as starting nodes. ",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "This class does the actual crawling. It is given a list of seed URLs, and effectively conducts a BFS using the seeds as starting nodes. ", 2)

# --- 6) Expand the "stops when it hits the limit" sentence.
$d.Content.Find.Execute(
    "It also keeps track of the depth of the search, and stops when it hits the limit.",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "It also keeps track of the depth of the search, and stops adding links to scrape when the current page it is on has reached that limit.", 2)

# --- 7) "SearchTool" heading: merge "SearchTool" + ":" runs, dropping the spell-check split.
$t = "SearchTool:"
$d.Content.Find.Execute($t, $false, $false, $false, $false, $false, $true, 1, $false, $t, 2)

# --- 8) Move the "_GoBack" bookmark from its own empty paragraph (just before "SearchTool:")
#        into the final paragraph, and extend the "console..." sentence with
#        " in the order of highest page ranks" right before the bookmark.
$d.Bookmarks.Item("_GoBack").Delete()

$d.Content.Find.Execute(
    "Finally, it prints the sites that keyword match to the console, and saves a list of them to a text file.",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "Finally, it prints the sites that keyword match to the console in the order of highest page ranks, and saves a list of them to a text file.", 2)

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$anchorRange = $lastPara.Range
$anchorText = "Finally, it prints the sites that keyword match to the console in the order of highest page ranks"
$anchorRange.Find.Execute($anchorText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$bmPos = $anchorRange.End
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

Write-Host "Done"
